$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1466.4166
$ws.Range("I40").Value = 1289.8
$ws.Range("K40").Value = 1289.8
$ws.Range("M40").Value = -1114.8
$ws.Range("H80").Value = 395.6875
$ws.Range("I80").Value = 194.83333
$ws.Range("J80").Value = 516.2
$ws.Range("K80").Value = 584.49999
$ws.Range("L80").Value = 1548.6
$ws.Range("M80").Value = 413.50001
$ws.Range("N80").Value = -3544.6
$ws.Range("H83").Value = 395.6875
$ws.Range("I83").Value = 194.83333
$ws.Range("J83").Value = 516.2
$ws.Range("K83").Value = 1753.49997
$ws.Range("L83").Value = 4645.8
$ws.Range("M83").Value = 3238.50003
$ws.Range("N83").Value = -14629.8
$ws.Range("H86").Value = 3488.4
$ws.Range("I86").Value = 3000
$ws.Range("J86").Value = 3610.5
$ws.Range("K86").Value = 3000
$ws.Range("L86").Value = 3610.5
$ws.Range("M86").Value = -1877
$ws.Range("N86").Value = -5856.5
$ws.Range("H89").Value = 3488.4
$ws.Range("I89").Value = 3000
$ws.Range("J89").Value = 3610.5
$ws.Range("K89").Value = 15000
$ws.Range("L89").Value = 18052.5
$ws.Range("M89").Value = -9384
$ws.Range("N89").Value = -29284.5
$ws.Range("H98").Value = 2264.9524
$ws.Range("I98").Value = 2264.9524
$ws.Range("K98").Value = 2264.9524
$ws.Range("M98").Value = -766.9524000000001
$ws.Range("H107").Value = 586.6316
$ws.Range("I107").Value = 519.625
$ws.Range("J107").Value = 944
$ws.Range("K107").Value = 519.625
$ws.Range("L107").Value = 944
$ws.Range("M107").Value = 1400.375
$ws.Range("N107").Value = -4784
$ws.Range("H116").Value = 3937
$ws.Range("J116").Value = 3875
$ws.Range("L116").Value = 3875
$ws.Range("N116").Value = -10759
$ws.Range("H122").Value = 2264.9524
$ws.Range("I122").Value = 2264.9524
$ws.Range("K122").Value = 6794.8572
$ws.Range("M122").Value = -4344.8572
$ws.Range("H130").Value = 105825.4
$ws.Range("J130").Value = 105825.4
$ws.Range("L130").Value = 105825.4
$ws.Range("N130").Value = -115865.4
$ws.Range("H132").Value = 1471.0769
$ws.Range("I132").Value = 1205.125
$ws.Range("K132").Value = 3615.375
$ws.Range("M132").Value = -1085.375
$ws.Range("H138").Value = 5954860.5
$ws.Range("J138").Value = 8931669
$ws.Range("L138").Value = 26795007
$ws.Range("N138").Value = -26805287

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9896.772000000001
$ws.Range("I32").Value = 5983.2705
$ws.Range("K32").Value = 5983.2705
$ws.Range("M32").Value = -5696.2705
$ws.Range("H61").Value = 3579.6191
$ws.Range("I61").Value = 2497.6099
$ws.Range("J61").Value = 5596.091
$ws.Range("K61").Value = 2497.6099
$ws.Range("L61").Value = 5596.091
$ws.Range("M61").Value = -2285.6099
$ws.Range("N61").Value = -6020.091
$ws.Range("H122").Value = 1046.2941
$ws.Range("I122").Value = 712.38464
$ws.Range("K122").Value = 2137.15392
$ws.Range("M122").Value = 312.8460800000003
$ws.Range("H132").Value = 2927.3572
$ws.Range("I132").Value = 2828.0244
$ws.Range("K132").Value = 8484.073199999999
$ws.Range("M132").Value = -5954.073199999999
$ws.Range("H136").Value = 3579.6191
$ws.Range("I136").Value = 2497.6099
$ws.Range("J136").Value = 5596.091
$ws.Range("K136").Value = 7492.8297
$ws.Range("L136").Value = 16788.273
$ws.Range("M136").Value = -4942.8297
$ws.Range("N136").Value = -21888.273

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2012.386
$ws.Range("I134").Value = 2022.2037
$ws.Range("J134").Value = 1835.6666
$ws.Range("K134").Value = 6066.6111
$ws.Range("L134").Value = 5506.9998
$ws.Range("M134").Value = -3531.6111
$ws.Range("N134").Value = -10576.9998

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 12500
$ws.Range("J50").Value = 12500
$ws.Range("L50").Value = 12500
$ws.Range("N50").Value = -13750
$ws.Range("H60").Value = 12833.5
$ws.Range("J60").Value = 11400.4
$ws.Range("L60").Value = 11400.4
$ws.Range("N60").Value = -12422.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 1031.0769
$ws.Range("I50").Value = 267.22223
$ws.Range("J50").Value = 2749.75
$ws.Range("K50").Value = 801.66669
$ws.Range("L50").Value = 8249.25
$ws.Range("M50").Value = -320.66669
$ws.Range("N50").Value = -9211.25
$ws.Range("H53").Value = 1031.0769
$ws.Range("I53").Value = 267.22223
$ws.Range("J53").Value = 2749.75
$ws.Range("K53").Value = 801.66669
$ws.Range("L53").Value = 8249.25
$ws.Range("M53").Value = -320.66669
$ws.Range("N53").Value = -9211.25
$ws.Range("H109").Value = 142.33333
$ws.Range("I109").Value = 142.33333
$ws.Range("K109").Value = 426.99999
$ws.Range("M109").Value = 613.00001
$ws.Range("H131").Value = 113218
$ws.Range("I131").Value = 500750
$ws.Range("J131").Value = 2494.5715
$ws.Range("K131").Value = 1502250
$ws.Range("L131").Value = 7483.7145
$ws.Range("M131").Value = -1497210
$ws.Range("N131").Value = -17563.7145
$ws.Range("H132").Value = 1449.25
$ws.Range("I132").Value = 1380.625
$ws.Range("J132").Value = 1723.75
$ws.Range("K132").Value = 12425.625
$ws.Range("L132").Value = 15513.75
$ws.Range("M132").Value = -9895.625
$ws.Range("N132").Value = -20573.75
$ws.Range("H137").Value = 1658.5555
$ws.Range("I137").Value = 1373
$ws.Range("J137").Value = 1923.7142
$ws.Range("K137").Value = 4119
$ws.Range("L137").Value = 5771.142599999999
$ws.Range("M137").Value = 981
$ws.Range("N137").Value = -15971.1426
$ws.Range("H140").Value = 2047.0667
$ws.Range("I140").Value = 1669.6923
$ws.Range("K140").Value = 5009.0769
$ws.Range("M140").Value = 170.9231

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 25000
$ws.Range("J26").Value = 25000
$ws.Range("L26").Value = 25000
$ws.Range("N26").Value = -25560
$ws.Range("H50").Value = 25000
$ws.Range("J50").Value = 25000
$ws.Range("L50").Value = 25000
$ws.Range("N50").Value = -25996
$ws.Range("H132").Value = 3733.1538
$ws.Range("I132").Value = 2446.4443
$ws.Range("K132").Value = 7339.3329
$ws.Range("M132").Value = -4809.3329
$ws.Range("H141").Value = 111353.5
$ws.Range("J141").Value = 111353.5
$ws.Range("L141").Value = 111353.5
$ws.Range("N141").Value = -121713.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4129.3076
$ws.Range("I40").Value = 2868.2
$ws.Range("J40").Value = 8333
$ws.Range("K40").Value = 2868.2
$ws.Range("L40").Value = 8333
$ws.Range("M40").Value = -2732.2
$ws.Range("N40").Value = -8605
$ws.Range("H43").Value = 9999
$ws.Range("J43").Value = 9999
$ws.Range("L43").Value = 9999
$ws.Range("N43").Value = -10385
$ws.Range("H132").Value = 3893.2354
$ws.Range("I132").Value = 3574.0625
$ws.Range("J132").Value = 9000
$ws.Range("K132").Value = 10722.1875
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -8192.1875
$ws.Range("N132").Value = -32060
$ws.Range("H136").Value = 3708.0605
$ws.Range("I136").Value = 3345.5334
$ws.Range("K136").Value = 10036.6002
$ws.Range("M136").Value = -7486.600199999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2158.2354
$ws.Range("I113").Value = 1885.5
$ws.Range("J113").Value = 2812.8
$ws.Range("K113").Value = 5656.5
$ws.Range("L113").Value = 8438.400000000001
$ws.Range("M113").Value = -3486.5
$ws.Range("N113").Value = -12778.4
$ws.Range("H122").Value = 1583.4839
$ws.Range("I122").Value = 1243.88
$ws.Range("J122").Value = 2998.5
$ws.Range("K122").Value = 3731.64
$ws.Range("L122").Value = 8995.5
$ws.Range("M122").Value = -1281.64
$ws.Range("N122").Value = -13895.5
$ws.Range("H125").Value = 29000
$ws.Range("J125").Value = 29000
$ws.Range("L125").Value = 29000
$ws.Range("N125").Value = -38840
$ws.Range("H136").Value = 2001.8108
$ws.Range("I136").Value = 1988.4482
$ws.Range("K136").Value = 5965.3446
$ws.Range("M136").Value = -3415.3446
